# Update "Lương" sheet: remove the now-unused "Phụ cấp" line items for
# CẦN THƠ and SÓC TRĂNG (code for phu cap tai cac co so khac), shifting
# all subsequent rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Row 3 = "Phụ cấp tại CẦN THƠ"
$ws.Rows.Item(3).Delete()

# After the first deletion everything shifted up by one row, so
# "Phụ cấp tại SÓC TRĂNG" (originally row 24) is now row 23.
$ws.Rows.Item(23).Delete()
